$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.565.87'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = '1.876.84'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('E4').Value = '  +0.99%  '
$ws.Range('D5').Value = '''313.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').Value = '''1.015'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('D7').Value = '''0.4801'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('D8').Value = '''0.3792'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.10%  '
$ws.Range('D9').Value = '''0.07393'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('D10').Value = '''0.9416'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.98%  '
$ws.Range('D11').Value = '''20.72'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.70%  '
$ws.Range('D12').Value = '''0.07875'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.41%  '
$ws.Range('D13').Value = '1.898.45'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('E14').Value = '  +2.85%  '
$ws.Range('D15').Value = '''6.609'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.40%  '
$ws.Range('D16').Value = '''91.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.36%  '
$ws.Range('D17').Value = '''1.017'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '''0.000008994'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.90%  '
$ws.Range('D19').Value = '''1.014'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('D20').Value = '''14.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.89%  '
$ws.Range('D21').Value = '27.592.62'
$ws.Range('E21').Value = '  +2.44%  '
$ws.Range('D22').Value = '''5.148'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.28%  '
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('D24').Value = '''1.968'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '''153.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = '''18.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.46%  '
$ws.Range('D27').Value = '''2.027'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('D28').Value = '''116.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('D29').Value = '''5.021'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').Value = '''0.08940'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('D31').Value = '''3.328'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').Value = '''1.215'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.19%  '
$ws.Range('E33').Value = '  +3.01%  '
$ws.Range('D34').Value = '''0.7527'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').Value = '''2.697'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = '''0.02080'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.75%  '
$ws.Range('D37').Value = '''1.123'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('D38').Value = '''0.05317'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('D40').Value = '''0.5382'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.45%  '
$ws.Range('E41').Value = '  +3.19%  '
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').Value = '''8.450'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.07%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''10.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.41%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.4856'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').Value = '''1.666'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = '''103.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('D49').Value = '''67.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.00%  '
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').Value = '''0.9028'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.23%  '
